$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (row 13) with 2021 data, following the pattern of the
# preceding row 12 (bold/centered/bordered style in column A, plain
# number cells elsewhere).
$row = 13
$prevRow = $row - 1

# Copy the formatting of the year-label cell (column A) from the row above
# so the new cell reuses the existing style (bold, centered, bordered).
$ws.Cells.Item($prevRow, 1).Copy($ws.Cells.Item($row, 1))
$ws.Cells.Item($row, 1).Value = "2021年"

$values = @{
    2  = 2496
    3  = 2520
    4  = 1142.05
    5  = 3918.28
    6  = 7686
    7  = 5935.96
    8  = 13.5466
    9  = 22635.34
    10 = 38
    11 = 8535
    12 = 1418.6876
    13 = 2912
    14 = 20165
    15 = 8.9092
}

foreach ($col in $values.Keys) {
    $ws.Cells.Item($row, $col).Value = $values[$col]
}
